$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.886.79"
$ws.Range("E2").Value = "  -3.26%  "
$ws.Range("D3").Value = "3.232.50"
$ws.Range("E3").Value = "  -3.78%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'538.75"
$ws.Range("D6").Value = "'136.50"
$ws.Range("E6").Value = "  -8.19%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "3.230.93"
$ws.Range("E8").Value = "  -3.85%  "
$ws.Range("E9").Value = "  -4.33%  "
$ws.Range("D10").Value = "'7.61"
$ws.Range("E10").Value = "  -4.01%  "
$ws.Range("E11").Value = "  -5.73%  "
$ws.Range("E12").Value = "  -4.29%  "
$ws.Range("D13").Value = "3.788.10"
$ws.Range("E13").Value = "  -3.79%  "
$ws.Range("E14").Value = "  -1.14%  "
$ws.Range("E15").Value = "  -7.09%  "
$ws.Range("D16").Value = "3.237.86"
$ws.Range("E16").Value = "  -4.04%  "
$ws.Range("E17").Value = "  -5.96%  "
$ws.Range("D18").Value = "58.905.01"
$ws.Range("E18").Value = "  -3.39%  "
$ws.Range("D19").Value = "'5.91"
$ws.Range("E19").Value = "  -7.26%  "
$ws.Range("D20").Value = "'13.29"
$ws.Range("E20").Value = "  -6.33%  "
$ws.Range("D21").Value = "'8.29"
$ws.Range("E21").Value = "  -6.35%  "
$ws.Range("D22").Value = "'361.47"
$ws.Range("E22").Value = "  -3.28%  "
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").Value = "'70.54"
$ws.Range("E24").Value = "  -6.54%  "
$ws.Range("E25").Value = "  -7.07%  "
$ws.Range("D26").Value = "3.373.93"
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").Value = "0.0₃0974"
$ws.Range("E27").Value = "  -10.44%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "'0.170"
$ws.Range("E28").Value = "  -3.51%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  -4.09%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").Value = "'1.94"
$ws.Range("E32").Value = "  -6.74%  "
$ws.Range("D33").Value = "'7.08"
$ws.Range("E33").Value = "  -7.54%  "
$ws.Range("D34").Value = "'21.99"
$ws.Range("E34").Value = "  -3.85%  "
$ws.Range("E35").Value = "  -4.67%  "
$ws.Range("D36").Value = "'163.10"
$ws.Range("E36").Value = "  -3.29%  "
$ws.Range("E37").Value = "  -7.85%  "
$ws.Range("D38").Value = "'6.42"
$ws.Range("E38").Value = "  -5.04%  "
$ws.Range("E39").Value = "  -6.81%  "
$ws.Range("D40").Value = "'26.49"
$ws.Range("E40").Value = "  -9.29%  "
$ws.Range("D41").Value = "'0.0711"
$ws.Range("E41").Value = "  -5.35%  "
$ws.Range("D42").Value = "3.264.42"
$ws.Range("E42").Value = "  -3.86%  "
$ws.Range("D43").Value = "'41.14"
$ws.Range("E43").Value = "  -2.70%  "
$ws.Range("D44").Value = "'0.716"
$ws.Range("E44").Value = "  -5.95%  "
$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D45").Value = "'1.10"
$ws.Range("E45").Value = "  -3.58%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value = "'4.03"
$ws.Range("E46").Value = "  -6.04%  "
$ws.Range("E47").Value = "  -6.38%  "
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("D49").Value = "2.297.64"
$ws.Range("E49").Value = "  -8.03%  "
$ws.Range("D50").Value = "'6.30"
$ws.Range("E50").Value = "  -5.69%  "
$ws.Range("E51").Value = "  -7.77%  "
